$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price / 1h-volume-change data and re-order a few
# rows whose relative ranking changed, as published by the
# "Updated cryptos list" GitHub Actions job.
#
# Column D ("Price") values are stored as plain text (they use the
# European "thousand.thousand.decimal" grouping, e.g. "2.310.76", so they
# can never be real numbers, but short values like "25.20" parse fine as a
# float). Force a text NumberFormat before assigning so COM/Excel does not
# silently coerce those values to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.675.93'
$ws.Range("E2").Value = '  -1.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.304.51'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.02'
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.23'
$ws.Range("E6").Value = '  -4.06%  '
$ws.Range("E7").Value = '  -5.17%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -5.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.43'
$ws.Range("E10").Value = '  -5.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.80'
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0789'
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.666.60'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.68'
$ws.Range("E16").Value = '  +4.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.311.47'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.620.46'
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.11'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.54'
$ws.Range("E22").Value = '  -5.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.91'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.99'
$ws.Range("E24").Value = '  -3.07%  '
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.53'
$ws.Range("E26").Value = '  -3.50%  '
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.96'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.20'
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.61'
$ws.Range("E30").Value = '  -6.43%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").Value = '  -9.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.19'
$ws.Range("E32").Value = '  -4.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.37'
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.62'
$ws.Range("E36").Value = '  +3.81%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.44'
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.93'
$ws.Range("E39").Value = '  -7.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  -5.25%  '
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("E42").Value = '  -5.10%  '
$ws.Range("E43").Value = '  -3.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.46'
$ws.Range("E44").Value = '  -7.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.991.67'
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.85'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0281'
$ws.Range("E47").Value = '  -4.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.23'
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  -5.73%  '
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.87'
$ws.Range("E51").Value = '  -2.44%  '
